# Updates as of 8th April 2020
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a handful of "Travelled From" entries that dropped their counts ---
$ws.Range("D7").Value  = "Spain[2] , Dubai"
$ws.Range("D20").Value = "Congolese"
$ws.Range("D21").Value = "Somali"
$ws.Range("D22").Value = "Pakistan,Kenyan(3)"
$ws.Range("D23").Value = "Nigerian"

# --- Append the two new daily rows (27 = 9-Apr-2020, 28 = 10-Apr-2020) ---
# Clone formatting from row 26 so the new rows match (date style on col A, etc.)
$ws.Rows.Item(26).Copy() | Out-Null
$ws.Rows.Item(27).PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(26).Copy() | Out-Null
$ws.Rows.Item(28).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 27
$ws.Cells.Item(27, 1).Value  = 43930
$ws.Cells.Item(27, 2).Value  = 5
$ws.Cells.Item(27, 3).Value  = 308
$ws.Cells.Item(27, 4).Value  = "Tanzania, UAE, UK"
$ws.Cells.Item(27, 5).Value  = "Nairobi (3) , Mombasa(2)"
$ws.Cells.Item(27, 6).Value  = 184
$ws.Cells.Item(27, 7).Value  = "Community(2),Import(3)"
$ws.Cells.Item(27, 8).Value  = 4
$ws.Cells.Item(27, 9).Value  = 0
$ws.Cells.Item(27, 12).Value = "39-77"
$ws.Cells.Item(27, 15).Value = 2
$ws.Cells.Item(27, 16).Value = 3

# Row 28
$ws.Cells.Item(28, 1).Value  = 43931
$ws.Cells.Item(28, 2).Value  = 5
$ws.Cells.Item(28, 3).Value  = 504
$ws.Cells.Item(28, 4).Value  = "UAE"
$ws.Cells.Item(28, 5).Value  = "Nairobi,Mombasa,Nyandarua(3)"
$ws.Cells.Item(28, 6).Value  = 189
$ws.Cells.Item(28, 7).Value  = "Community(4), Imported"
$ws.Cells.Item(28, 8).Value  = 10
$ws.Cells.Item(28, 9).Value  = 0
$ws.Cells.Item(28, 12).Value = "27-58"
$ws.Cells.Item(28, 15).Value = 3
$ws.Cells.Item(28, 16).Value = 2

# --- Scroll / selection bookkeeping to match the author's saved view ---
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G26").Select() | Out-Null
